$d = $word.ActiveDocument

$targets = @(
    "Folders with more in them are rendered as visibly larger wormholes",
    "Will have to implement a secondary traversal algorithm"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    foreach ($t in $targets) {
        if ($text -like "*$t*") {
            $p.Range.Font.StrikeThrough = $true
        }
    }
}
